{"js": "// Translate the English labels/cells of the \"Facilitators guidelines -\n// Moebius\" table into French. Each entry below is an exact, whole-run\n// source string paired with its French replacement (mirrors the XML diff,\n// one `w:t` run per hit; \"Assist the process, provoke thoughts\" repeats\n// 6x across the document and every occurrence gets the same translation).\nconst replacements = [\n  [\"Video Title\", \"Titre de la vid\u00e9o\"],\n  [\"Topic\", \"Rubrique\"],\n  [\"Geometry\", \"G\u00e9om\u00e9trie\"],\n  [\"Aim(s)\", \"Objectif(s)\"],\n  [\"Length\", \"Dur\u00e9e\"],\n  [\"Camp Location\", \"Lieu du camp\"],\n  [\"Facilitators\", \"Animateurs\"],\n  [\"N. of students\", \"N. des \u00e9tudiants\"],\n  [\"Resources\", \"Les ressources\"],\n  [\"needed\", \"n\u00e9cessaires\"],\n  [\"Preparations\", \"Pr\u00e9parations\"],\n  [\"Video time\", \"Temps de la vid\u00e9o\"],\n  [\"What facilitator does\", \"Ce que fait le facilitateur\"],\n  [\"What learners do\", \"Ce que font les apprenants\"],\n  [\"General VMC Video Introduction\", \"Vid\u00e9o g\u00e9n\u00e9rale introduisant le CVM\"],\n  [\"Video Introduction\", \"Video d'introduction\"],\n  [\"Introduction of the first experiment\", \"Introduction de la premi\u00e8re exp\u00e9rimentation\"],\n  [\"Assist the process, provoke thoughts\", \"Faciliter le processus, susciter des pens\u00e9es\"],\n];\n\nconst body = context.document.body;\n\nfor (const [search, replacement] of replacements) {\n  const found = body.search(search, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.load(\"text\");\n  }\n  await context.sync();\n\n  for (const range of found.items) {\n    // Only swap ranges whose text is an exact match for the whole target\n    // string (guards against a search hit landing inside a longer run).\n    if (range.text === search) {\n      range.insertText(replacement, \"Replace\");\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Translate the English labels/cells of the \"Facilitators guidelines -\n# Moebius\" table into French. Each pair below is an exact, whole-run\n# source string matched with Find.Execute/wdReplaceAll and swapped for its\n# French translation (mirrors the XML diff; \"Assist the process, provoke\n# thoughts\" repeats 6x across the document and every hit gets replaced).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"Video Title\", \"Titre de la vid\u00e9o\"),\n    @(\"Topic\", \"Rubrique\"),\n    @(\"Geometry\", \"G\u00e9om\u00e9trie\"),\n    @(\"Aim(s)\", \"Objectif(s)\"),\n    @(\"Length\", \"Dur\u00e9e\"),\n    @(\"Camp Location\", \"Lieu du camp\"),\n    @(\"Facilitators\", \"Animateurs\"),\n    @(\"N. of students\", \"N. des \u00e9tudiants\"),\n    @(\"Resources\", \"Les ressources\"),\n    @(\"needed\", \"n\u00e9cessaires\"),\n    @(\"Preparations\", \"Pr\u00e9parations\"),\n    @(\"Video time\", \"Temps de la vid\u00e9o\"),\n    @(\"What facilitator does\", \"Ce que fait le facilitateur\"),\n    @(\"What learners do\", \"Ce que font les apprenants\"),\n    @(\"General VMC Video Introduction\", \"Vid\u00e9o g\u00e9n\u00e9rale introduisant le CVM\"),\n    @(\"Video Introduction\", \"Video d'introduction\"),\n    @(\"Introduction of the first experiment\", \"Introduction de la premi\u00e8re exp\u00e9rimentation\"),\n    @(\"Assist the process, provoke thoughts\", \"Faciliter le processus, susciter des pens\u00e9es\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nWrite-Output \"done\"\n"}
